# Updates the cryptos price/volume table with freshly scraped values.
# Numeric-looking "D" (Price) values are written with NumberFormat "@"
# (Text) first so Excel doesn't coerce strings like "7.00" / "12.00" into
# numbers (which would drop the meaningful trailing zeros / formatting);
# the style is reset back to "Normal" right after so no stray cell style
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.364.98'
$ws.Range('E2').Value = '  +0.90%  '

$ws.Range('D3').Value = '2.366.96'
$ws.Range('E3').Value = '  +2.94%  '

$ws.Range('E4').Value = '  -0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.79'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '104.74'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.15%  '

$ws.Range('E7').Value = '  -1.36%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.41%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.29'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.26%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.03'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.96%  '

$ws.Range('E12').Value = '  -0.61%  '

$ws.Range('E13').Value = '  -0.58%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.00'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.77%  '

$ws.Range('D15').Value = '2.734.23'
$ws.Range('E15').Value = '  +2.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.62'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.36%  '

$ws.Range('D17').Value = '2.378.01'
$ws.Range('E17').Value = '  +3.42%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.818'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.06%  '

$ws.Range('D19').Value = '43.309.73'
$ws.Range('E19').Value = '  +0.78%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.10%  '

$ws.Range('E21').Value = '  +0.61%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.47%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.43'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.56%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '242.19'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.02%  '

$ws.Range('E25').Value = '  +2.21%  '

$ws.Range('E26').Value = '  +0.69%  '

$ws.Range('E27').Value = '  +0.72%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '26.00'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +7.34%  '

$ws.Range('E29').Value = '  +8.85%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.92'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.85%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.13%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '161.72'
$ws.Range('D32').Style = 'Normal'

$ws.Range('E33').Value = '  -0.48%  '

$ws.Range('E34').Value = '  -0.16%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.37'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.74%  '

$ws.Range('E36').Value = '  +6.48%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.12'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.94%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.68'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +11.98%  '

$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.95'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +5.95%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0739'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.23%  '

$ws.Range('E41').Value = '  +0.70%  '

$ws.Range('E42').Value = '  -0.83%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.97%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.23'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +5.32%  '

$ws.Range('D45').Value = '2.004.21'
$ws.Range('E45').Value = '  +1.92%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0291'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.69%  '

$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.79%  '

$ws.Range('E48').Value = '  +6.16%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '58.05'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.74%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.93'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.42%  '

$ws.Range('D51').Value = '2.574.56'
$ws.Range('E51').Value = '  +1.89%  '

